$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Foglio1")

# Sprint 2 burndown: points completed on 2019-12-09 (row 40) corrected from 7 to 15.
# Columns F (cumulative completed) and G (points left) are formulas that
# recalculate automatically for the dependent rows (40-47).
$ws.Range("B40").Value = 15

# Restore the on-screen selection/scroll state.
$ws.Range("A37").Select()
$ws.Range("B41").Select()
